# Add two new columns ("Aportes" and "Turno") to the sales sheet, fin clase 3 in
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (E1, F1) ---
$ws.Range("E1").Value = "Aportes"
$ws.Range("F1").Value = "Turno"

# Copy the existing header style (used by B1:D1) onto the new headers
$ws.Range("B1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data values (E2:F11) ---
$aportes = @(
    0.1139778962782939,
    0.04743262275842433,
    0.1976395354900047,
    0.1655148352056989,
    0.1026941048696355,
    0.05807755618836215,
    0.156111127963929,
    0.07690169178220757,
    0.1674846750856551,
    0.0372984697196711
)

$turno = @(
    "Nocturno",
    "Nocturno",
    "Diurno",
    "Nocturno",
    "Nocturno",
    "Diurno",
    "Diurno",
    "Diurno",
    "Nocturno",
    "Diurno"
)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $aportes[$i]
    $ws.Cells.Item($row, 6).Value = $turno[$i]
}
